$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("June")

$ws.Range("B2").Value = 1691
$ws.Range("C2").Value = 1104
$ws.Range("D2").Value = 587
$ws.Range("E2").Value = 'We borrowerd more than we lent'
$ws.Range("F2").Value = ""
$ws.Range("G2").Value = '1.53 : 1'

$ws.Range("B3").Value = 477
$ws.Range("C3").Value = 448
$ws.Range("D3").Value = 29
$ws.Range("E3").Value = 'We borrowerd more than we lent'
$ws.Range("F3").Value = ""
$ws.Range("G3").Value = '1.06 : 1'

$ws.Range("B4").Value = 1394
$ws.Range("C4").Value = 1350
$ws.Range("D4").Value = 44
$ws.Range("E4").Value = 'We borrowerd more than we lent'
$ws.Range("F4").Value = ""
$ws.Range("G4").Value = '1.03 : 1'

$ws.Range("B5").Value = 84
$ws.Range("C5").Value = 119
$ws.Range("D5").Value = -35
$ws.Range("E5").Value = ""
$ws.Range("F5").Value = 'We lent more than we borrowed'
$ws.Range("G5").Value = '0.71 : 1'

$ws.Range("B6").Value = 1047
$ws.Range("C6").Value = 1308
$ws.Range("D6").Value = -261
$ws.Range("E6").Value = ""
$ws.Range("F6").Value = 'We lent more than we borrowed'
$ws.Range("G6").Value = '0.80 : 1'

$ws.Range("B7").Value = 174
$ws.Range("C7").Value = 221
$ws.Range("D7").Value = -47
$ws.Range("E7").Value = ""
$ws.Range("F7").Value = 'We lent more than we borrowed'
$ws.Range("G7").Value = '0.79 : 1'

$ws.Range("B8").Value = 127
$ws.Range("C8").Value = 242
$ws.Range("D8").Value = -115
$ws.Range("E8").Value = ""
$ws.Range("F8").Value = 'We lent more than we borrowed'
$ws.Range("G8").Value = '0.52 : 1'

$ws.Range("B9").Value = 46
$ws.Range("C9").Value = 55
$ws.Range("D9").Value = -9
$ws.Range("E9").Value = ""
$ws.Range("F9").Value = 'We lent more than we borrowed'
$ws.Range("G9").Value = '0.84 : 1'

$ws.Range("B10").Value = 3
$ws.Range("C10").Value = 13
$ws.Range("D10").Value = -10
$ws.Range("E10").Value = ""
$ws.Range("F10").Value = 'We lent more than we borrowed'
$ws.Range("G10").Value = '0.23 : 1'

$ws.Range("B11").Value = 0
$ws.Range("C11").Value = 0
$ws.Range("D11").Value = 0
$ws.Range("E11").Value = ""
$ws.Range("F11").Value = ""
$ws.Range("G11").Value = ""

$ws.Range("B12").Value = 16
$ws.Range("C12").Value = 24
$ws.Range("D12").Value = -8
$ws.Range("E12").Value = ""
$ws.Range("F12").Value = 'We lent more than we borrowed'
$ws.Range("G12").Value = '0.67 : 1'

$ws.Range("B13").Value = 138
$ws.Range("C13").Value = 95
$ws.Range("D13").Value = 43
$ws.Range("E13").Value = 'We borrowerd more than we lent'
$ws.Range("F13").Value = ""
$ws.Range("G13").Value = '1.45 : 1'

$ws.Range("B14").Value = 137
$ws.Range("C14").Value = 290
$ws.Range("D14").Value = -153
$ws.Range("E14").Value = ""
$ws.Range("F14").Value = 'We lent more than we borrowed'
$ws.Range("G14").Value = '0.47 : 1'

$ws.Range("B15").Value = 81
$ws.Range("C15").Value = 97
$ws.Range("D15").Value = -16
$ws.Range("E15").Value = ""
$ws.Range("F15").Value = 'We lent more than we borrowed'
$ws.Range("G15").Value = '0.84 : 1'

$ws.Range("B16").Value = 26
$ws.Range("C16").Value = 144
$ws.Range("D16").Value = -118
$ws.Range("E16").Value = ""
$ws.Range("F16").Value = 'We lent more than we borrowed'
$ws.Range("G16").Value = '0.18 : 1'

$ws.Range("B17").Value = 637
$ws.Range("C17").Value = 457
$ws.Range("D17").Value = 180
$ws.Range("E17").Value = 'We borrowerd more than we lent'
$ws.Range("F17").Value = ""
$ws.Range("G17").Value = '1.39 : 1'

$ws.Range("B18").Value = 19
$ws.Range("C18").Value = 106
$ws.Range("D18").Value = -87
$ws.Range("E18").Value = ""
$ws.Range("F18").Value = 'We lent more than we borrowed'
$ws.Range("G18").Value = '0.18 : 1'

$ws.Range("B19").Value = 590
$ws.Range("C19").Value = 417
$ws.Range("D19").Value = 173
$ws.Range("E19").Value = 'We borrowerd more than we lent'
$ws.Range("F19").Value = ""
$ws.Range("G19").Value = '1.41 : 1'

$ws.Range("B20").Value = 1
$ws.Range("C20").Value = 61
$ws.Range("D20").Value = -60
$ws.Range("E20").Value = ""
$ws.Range("F20").Value = 'We lent more than we borrowed'
$ws.Range("G20").Value = '0.02 : 1'

$ws.Range("B21").Value = 536
$ws.Range("C21").Value = 402
$ws.Range("D21").Value = 134
$ws.Range("E21").Value = 'We borrowerd more than we lent'
$ws.Range("F21").Value = ""
$ws.Range("G21").Value = '1.33 : 1'

$ws.Range("B22").Value = 28
$ws.Range("C22").Value = 65
$ws.Range("D22").Value = -37
$ws.Range("E22").Value = ""
$ws.Range("F22").Value = 'We lent more than we borrowed'
$ws.Range("G22").Value = '0.43 : 1'

$ws.Range("B23").Value = 676
$ws.Range("C23").Value = 458
$ws.Range("D23").Value = 218
$ws.Range("E23").Value = 'We borrowerd more than we lent'
$ws.Range("F23").Value = ""
$ws.Range("G23").Value = '1.48 : 1'

$ws.Range("B24").Value = 1731
$ws.Range("C24").Value = 1389
$ws.Range("D24").Value = 342
$ws.Range("E24").Value = 'We borrowerd more than we lent'
$ws.Range("F24").Value = ""
$ws.Range("G24").Value = '1.25 : 1'

$ws.Range("B25").Value = 160
$ws.Range("C25").Value = 417
$ws.Range("D25").Value = -257
$ws.Range("E25").Value = ""
$ws.Range("F25").Value = 'We lent more than we borrowed'
$ws.Range("G25").Value = '0.38 : 1'

$ws.Range("B26").Value = 0
$ws.Range("C26").Value = 0
$ws.Range("D26").Value = 0
$ws.Range("E26").Value = ""
$ws.Range("F26").Value = ""
$ws.Range("G26").Value = ""

$ws.Range("B27").Value = 245
$ws.Range("C27").Value = 213
$ws.Range("D27").Value = 32
$ws.Range("E27").Value = 'We borrowerd more than we lent'
$ws.Range("F27").Value = ""
$ws.Range("G27").Value = '1.15 : 1'

$ws.Range("B28").Value = 60
$ws.Range("C28").Value = 78
$ws.Range("D28").Value = -18
$ws.Range("E28").Value = ""
$ws.Range("F28").Value = 'We lent more than we borrowed'
$ws.Range("G28").Value = '0.77 : 1'

$ws.Range("B29").Value = 595
$ws.Range("C29").Value = 456
$ws.Range("D29").Value = 139
$ws.Range("E29").Value = 'We borrowerd more than we lent'
$ws.Range("F29").Value = ""
$ws.Range("G29").Value = '1.30 : 1'

$ws.Range("B30").Value = 29
$ws.Range("C30").Value = 48
$ws.Range("D30").Value = -19
$ws.Range("E30").Value = ""
$ws.Range("F30").Value = 'We lent more than we borrowed'
$ws.Range("G30").Value = '0.60 : 1'

$ws.Range("B31").Value = 62
$ws.Range("C31").Value = 299
$ws.Range("D31").Value = -237
$ws.Range("E31").Value = ""
$ws.Range("F31").Value = 'We lent more than we borrowed'
$ws.Range("G31").Value = '0.21 : 1'

$ws.Range("B32").Value = 399
$ws.Range("C32").Value = 468
$ws.Range("D32").Value = -69
$ws.Range("E32").Value = ""
$ws.Range("F32").Value = 'We lent more than we borrowed'
$ws.Range("G32").Value = '0.85 : 1'

$ws.Range("B33").Value = 420
$ws.Range("C33").Value = 530
$ws.Range("D33").Value = -110
$ws.Range("E33").Value = ""
$ws.Range("F33").Value = 'We lent more than we borrowed'
$ws.Range("G33").Value = '0.79 : 1'

$ws.Range("B34").Value = 244
$ws.Range("C34").Value = 90
$ws.Range("D34").Value = 154
$ws.Range("E34").Value = 'We borrowerd more than we lent'
$ws.Range("F34").Value = ""
$ws.Range("G34").Value = '2.71 : 1'

$ws.Range("B35").Value = 996
$ws.Range("C35").Value = 975
$ws.Range("D35").Value = 21
$ws.Range("E35").Value = 'We borrowerd more than we lent'
$ws.Range("F35").Value = ""
$ws.Range("G35").Value = '1.02 : 1'

$ws.Range("B36").Value = 202
$ws.Range("C36").Value = 355
$ws.Range("D36").Value = -153
$ws.Range("E36").Value = ""
$ws.Range("F36").Value = 'We lent more than we borrowed'
$ws.Range("G36").Value = '0.57 : 1'

$ws.Range("B37").Value = 585
$ws.Range("C37").Value = 341
$ws.Range("D37").Value = 244
$ws.Range("E37").Value = 'We borrowerd more than we lent'
$ws.Range("F37").Value = ""
$ws.Range("G37").Value = '1.72 : 1'

$ws.Range("B38").Value = 15
$ws.Range("C38").Value = 199
$ws.Range("D38").Value = -184
$ws.Range("E38").Value = ""
$ws.Range("F38").Value = 'We lent more than we borrowed'
$ws.Range("G38").Value = '0.08 : 1'

$ws.Range("B39").Value = 0
$ws.Range("C39").Value = 10
$ws.Range("D39").Value = -10
$ws.Range("E39").Value = ""
$ws.Range("F39").Value = 'We lent more than we borrowed'
$ws.Range("G39").Value = '0.00 : 1'

$ws.Range("B40").Value = 0
$ws.Range("C40").Value = 22
$ws.Range("D40").Value = -22
$ws.Range("E40").Value = ""
$ws.Range("F40").Value = 'We lent more than we borrowed'
$ws.Range("G40").Value = '0.00 : 1'

$ws.Range("B41").Value = 0
$ws.Range("C41").Value = 9
$ws.Range("D41").Value = -9
$ws.Range("E41").Value = ""
$ws.Range("F41").Value = 'We lent more than we borrowed'
$ws.Range("G41").Value = '0.00 : 1'

$ws.Range("B42").Value = 0
$ws.Range("C42").Value = 1
$ws.Range("D42").Value = -1
$ws.Range("E42").Value = ""
$ws.Range("F42").Value = 'We lent more than we borrowed'
$ws.Range("G42").Value = '0.00 : 1'

$ws.Range("B43").Value = 0
$ws.Range("C43").Value = 0
$ws.Range("D43").Value = 0
$ws.Range("E43").Value = ""
$ws.Range("F43").Value = ""
$ws.Range("G43").Value = ""

$ws.Range("B44").Value = 191
$ws.Range("C44").Value = 89
$ws.Range("D44").Value = 102
$ws.Range("E44").Value = 'We borrowerd more than we lent'
$ws.Range("F44").Value = ""
$ws.Range("G44").Value = '2.15 : 1'

$ws.Range("B45").Value = 115
$ws.Range("C45").Value = 154
$ws.Range("D45").Value = -39
$ws.Range("E45").Value = ""
$ws.Range("F45").Value = 'We lent more than we borrowed'
$ws.Range("G45").Value = '0.75 : 1'

$ws.Range("B46").Value = 467
$ws.Range("C46").Value = 687
$ws.Range("D46").Value = -220
$ws.Range("E46").Value = ""
$ws.Range("F46").Value = 'We lent more than we borrowed'
$ws.Range("G46").Value = '0.68 : 1'

$ws.Range("B47").Value = 780
$ws.Range("C47").Value = 589
$ws.Range("D47").Value = 191
$ws.Range("E47").Value = 'We borrowerd more than we lent'
$ws.Range("F47").Value = ""
$ws.Range("G47").Value = '1.32 : 1'

$ws.Range("B48").Value = 282
$ws.Range("C48").Value = 706
$ws.Range("D48").Value = -424
$ws.Range("E48").Value = ""
$ws.Range("F48").Value = 'We lent more than we borrowed'
$ws.Range("G48").Value = '0.40 : 1'

$ws.Range("B49").Value = 309
$ws.Range("C49").Value = 213
$ws.Range("D49").Value = 96
$ws.Range("E49").Value = 'We borrowerd more than we lent'
$ws.Range("F49").Value = ""
$ws.Range("G49").Value = '1.45 : 1'

$ws.Range("B50").Value = 945
$ws.Range("C50").Value = 573
$ws.Range("D50").Value = 372
$ws.Range("E50").Value = 'We borrowerd more than we lent'
$ws.Range("F50").Value = ""
$ws.Range("G50").Value = '1.65 : 1'

$ws.Range("B51").Value = 174
$ws.Range("C51").Value = 106
$ws.Range("D51").Value = 68
$ws.Range("E51").Value = 'We borrowerd more than we lent'
$ws.Range("F51").Value = ""
$ws.Range("G51").Value = '1.64 : 1'

$ws.Range("B52").Value = 318
$ws.Range("C52").Value = 532
$ws.Range("D52").Value = -214
$ws.Range("E52").Value = ""
$ws.Range("F52").Value = 'We lent more than we borrowed'
$ws.Range("G52").Value = '0.60 : 1'

$ws.Range("B53").Value = 134
$ws.Range("C53").Value = 225
$ws.Range("D53").Value = -91
$ws.Range("E53").Value = ""
$ws.Range("F53").Value = 'We lent more than we borrowed'
$ws.Range("G53").Value = '0.60 : 1'

$ws.Range("B54").Value = 26
$ws.Range("C54").Value = 178
$ws.Range("D54").Value = -152
$ws.Range("E54").Value = ""
$ws.Range("F54").Value = 'We lent more than we borrowed'
$ws.Range("G54").Value = '0.15 : 1'

$ws.Range("B55").Value = 234
$ws.Range("C55").Value = 218
$ws.Range("D55").Value = 16
$ws.Range("E55").Value = 'We borrowerd more than we lent'
$ws.Range("F55").Value = ""
$ws.Range("G55").Value = '1.07 : 1'

